$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$caseQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.stage_of_disease IN ['Vb']
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# B2 (CasesTab query) drops its trailing "Cohort" line (and the now-dangling
# trailing comma on the "Response to Treatment" line above it).
$ws.Range("B2").Value = $caseQuery

# Row heights settle to a uniform wrapped-text height after the edit
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Restore selection/scroll state to B2 (top-left resets to A1 automatically)
$ws.Range("B2").Select()
